# Applies the "Updated cryptos list" GitHub Actions data refresh to Sheet1.
# For each changed row: Price (col D) and/or Volume(1h) (col E) inline strings are updated
# to the new scraped values. Column D values that look numeric (so Excel would otherwise
# auto-convert them to a Number and mangle formatting like trailing zeros) are forced back
# to Text via NumberFormat "@" before the write, matching the original inlineStr text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.155.91"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.478.65"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.17"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.23"
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "2.478.36"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.93"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "2.939.68"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.54"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "67.068.96"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000170"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "2.524.89"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.65"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.00"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.10"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.97"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.23"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.80"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  +4.03%  "
$ws.Range("D28").Value = "2.603.95"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "0.0₃0908"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "508.95"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.70"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.04"
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.71"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.13"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.69"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.82"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("E45").Value = "  +3.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.50"
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("D47").Value = "0.0₆0264"
$ws.Range("E47").Value = "  +4.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.48"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.515"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.57"
$ws.Range("E51").Value = "  -0.28%  "
